# Applies the "Saldo" export update:
#   1. Insert new row 004479287 / ANA / 100000 right before account 004482102
#   2. Remove row 005046790 / BEATRIZ / 20000
#   3. Rename account 004641487 (LAILA) to 004643737 / LARA (balance unchanged)
#      and drop the now-duplicate old 004643737 / LARA row beneath it
#   4. Insert new row 004586209 / ROBERIO / 6400 right before account 004425965
#   5. Change the balance of account 004425965 (CAROLLINA) from 5014.27 to 1319.14
#   6. Remove row 004550415 / DIOGO / 1650
#   7. Insert new row 004444164 / ANA / 839.84 right after account 005009026

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$searchCol = $ws.Range("A1:A300")

function Get-RowOf($account) {
    $hit = $searchCol.Find($account)
    return $hit.Row
}

function Insert-AccountRow($beforeAccount, $account, $name, $saldo) {
    $r = Get-RowOf $beforeAccount
    $ws.Rows.Item($r).Insert()
    $ws.Cells.Item($r, 1).Value = "'" + $account
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $saldo
}

function Insert-AccountRowAfter($afterAccount, $account, $name, $saldo) {
    $r = (Get-RowOf $afterAccount) + 1
    $ws.Rows.Item($r).Insert()
    $ws.Cells.Item($r, 1).Value = "'" + $account
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $saldo
}

function Remove-AccountRow($account) {
    $r = Get-RowOf $account
    $ws.Rows.Item($r).Delete()
}

# 1. New account at the very top of the data (before 004482102)
Insert-AccountRow "004482102" "004479287" "ANA" 100000

# 2. Drop the BEATRIZ / 005046790 row entirely
Remove-AccountRow "005046790"

# 3. 004641487 / LAILA becomes 004643737 / LARA, balance (9133.77) stays the same
$r = Get-RowOf "004641487"
$ws.Cells.Item($r, 1).Value = "'004643737"
$ws.Cells.Item($r, 2).Value = "LARA"
# the old second 004643737 / LARA / 9133.77 row is now a duplicate - remove it
Remove-AccountRow "004643737"

# 4. New account just before CAROLLINA (004425965)
Insert-AccountRow "004425965" "004586209" "ROBERIO" 6400

# 5. CAROLLINA's balance drops from 5014.27 to 1319.14
$r = Get-RowOf "004425965"
$ws.Cells.Item($r, 3).Value = 1319.14

# 6. Drop the DIOGO / 004550415 row entirely
Remove-AccountRow "004550415"

# 7. New account right after EDMUR (005009026)
Insert-AccountRowAfter "005009026" "004444164" "ANA" 839.84
